$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("101_1")

$ws.Range("B17").Value = 1
$ws.Range("B17").HorizontalAlignment = -4152

$ws.Range("B26").Value = 1
$ws.Range("B26").HorizontalAlignment = -4152
